# Updates the "API ĐIỀU KHIỂN" bullet list:
#  - Paragraphs describing LEFT / RIGHT / BACK commands were re-typed as
#    single consolidated runs (their visible text is unchanged, only the
#    run-splitting collapses), and a brand-new bullet describing the
#    #LINEXSPEEDY! command is appended after the "Chạy lùi theo line" item.

$d = $word.ActiveDocument

# --- 1) "Quay trái 90 độ" bullet: collapse its runs into one (text is
#        unchanged - it was already split across "...tốc độ. " / "Ví dụ: "
#        / "quay trái" / " với tốc độ ... #" / "LEFT" / "150!").
$t1 = "Quay trái 90 độ: #LEFTXXX! Với XXX là tốc độ. Ví dụ: quay trái với tốc độ 150mm/s lệnh điều khiển sẽ là: #LEFT150!"
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

# --- 2) "Quay phải 90 độ" bullet: same kind of run collapse.
$t2 = "Quay phải 90 độ: #RIGHTXXX! Với XXX là tốc độ. Ví dụ: quay phải với tốc độ 150mm/s lệnh điều khiển sẽ là: #LEFT150!"
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

# --- 3) "Quay 180 độ" bullet: same kind of run collapse.
$t3 = "Quay 180 độ: #BACKXXX! Với XXX là tốc độ. Ví dụ: quay 180 độ với tốc độ 150mm/s lệnh điều khiển sẽ là: #BACK150!"
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# --- 4) "Đi chậm và dừng hẳn" bullet is untouched by the diff (already a
#        single run) - nothing to do.

# --- 5) "Chạy lùi theo line" bullet: same kind of run collapse.
$t5 = "Chạy lùi theo line: #RUNBACKXXX! Với XXX là tốc độ. Ví dụ: chạy lùi theo line với tốc độ 150mm/s lệnh điều khiển sẽ là: #RUNBACK150!"
$d.Content.Find.Execute($t5, $true, $false, $false, $false, $false, $true, 1, $false, $t5, 2) | Out-Null

# --- 6) Brand-new bullet for the #LINEXSPEEDY! command, added right after
#        the "Chạy lùi theo line" bullet, same list (ListParagraph / numId 3).
# Locate the paragraph whose text is the "Chạy lùi theo line" bullet and
# insert a new paragraph right after it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd("`r") -eq $t5) {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs($i + 1)
        $newPara.Range.Text = "Chạy đếm đủ line rồi dừng: #LINEXSPEEDY! Với X là số LINE cần qua, Y là tốc độ để chạy quãng đường đó."
        break
    }
}
